# Add season record columns (Wins, Losses, Ties) to the roster sheet.
# The original sheet tracked team statistics per-player but did not
# include the team's overall season record. This adds three new
# columns (AD, AE, AF) with that information.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): add the three new column headers, matching the
# existing header formatting (bold, bordered, centered) by copying the
# style from the neighboring header cell.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows (2-46): every row gets the same season record values.
for ($r = 2; $r -le 46; $r++) {
    $ws.Range("AD$r").Value = 85
    $ws.Range("AE$r").Value = 77
    $ws.Range("AF$r").Value = 0
}

Write-Host "Season record columns added."
